# Update the Test Suite Statistics sheet with the latest test counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the SourceProviderRDO row (row 3) counts from 9 to 10.
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 10

# Recalculate so dependent formulas (G4, G5, etc.) pick up the new totals.
$excel.Calculate()

# Move the active selection to C4, as left by the editor after the update.
$ws.Activate()
$ws.Range("C4").Select()
